# Update latest output (run 174)
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Schedule": rows 3-5 get new optimisation numbers, and a new row 6
# is appended (dimension grows from A1:F5 to A1:F6).
# ---------------------------------------------------------------------------
$schedule = $wb.Worksheets.Item("Schedule")

# Row 3
$schedule.Cells.Item(3, 1).Value = 46066.33333333334
$schedule.Cells.Item(3, 3).Value = 8
$schedule.Cells.Item(3, 4).Value = 30.24
$schedule.Cells.Item(3, 5).Value = 923.0638170000002
$schedule.Cells.Item(3, 6).Value = 30.52459712301588

# Row 4
$schedule.Cells.Item(4, 1).Value = 46066.83333333334
$schedule.Cells.Item(4, 2).Value = 46067
$schedule.Cells.Item(4, 3).Value = 4
$schedule.Cells.Item(4, 4).Value = 15.12
$schedule.Cells.Item(4, 5).Value = 550.5185009999999
$schedule.Cells.Item(4, 6).Value = 36.40995376984127

# Row 5
$schedule.Cells.Item(5, 1).Value = 46067.04166666666
$schedule.Cells.Item(5, 2).Value = 46067.20833333334
$schedule.Cells.Item(5, 3).Value = 4
$schedule.Cells.Item(5, 4).Value = 15.12
$schedule.Cells.Item(5, 5).Value = 406.95022875
$schedule.Cells.Item(5, 6).Value = 26.9146976686508

# Row 6 (new) - copy the date/time display format from row 5's A/B cells
$schedule.Range("A5:B5").Copy()
$schedule.Range("A6:B6").PasteSpecial(-4122)

$schedule.Cells.Item(6, 1).Value = 46067.375
$schedule.Cells.Item(6, 2).Value = 46067.70833333334
$schedule.Cells.Item(6, 3).Value = 8
$schedule.Cells.Item(6, 4).Value = 30.24
$schedule.Cells.Item(6, 5).Value = 251.7862035
$schedule.Cells.Item(6, 6).Value = 8.326263343253968

# ---------------------------------------------------------------------------
# Sheet "Detailed": refreshed price forecast/historical values, pump-status
# flags, and forecast/historical classifications for the latest run.
# ---------------------------------------------------------------------------
$detailed = $wb.Worksheets.Item("Detailed")

$detailed.Cells.Item(17, 5).Value = "OFF"

$detailed.Cells.Item(38, 2).Value = 49.00286

$detailed.Cells.Item(39, 2).Value = 39.58239

$detailed.Cells.Item(40, 2).Value = 66.37485
$detailed.Cells.Item(40, 3).Value = "historical"

$detailed.Cells.Item(41, 2).Value = 73.19
$detailed.Cells.Item(41, 3).Value = "historical"

$detailed.Cells.Item(42, 2).Value = 67.05637
$detailed.Cells.Item(42, 3).Value = "historical"
$detailed.Cells.Item(42, 5).Value = "ON"

$detailed.Cells.Item(43, 2).Value = 71.80374999999999
$detailed.Cells.Item(43, 3).Value = "historical"

$detailed.Cells.Item(44, 2).Value = 68.15342
$detailed.Cells.Item(44, 3).Value = "historical"

$detailed.Cells.Item(45, 2).Value = 71.40000000000001
$detailed.Cells.Item(45, 3).Value = "historical"

$detailed.Cells.Item(46, 2).Value = 71.40000000000001
$detailed.Cells.Item(46, 3).Value = "historical"

$detailed.Cells.Item(47, 2).Value = 72.51486
$detailed.Cells.Item(47, 3).Value = "historical"

$detailed.Cells.Item(48, 2).Value = 70.80972
$detailed.Cells.Item(48, 3).Value = "historical"

$detailed.Cells.Item(49, 2).Value = 71.49624

$detailed.Cells.Item(50, 2).Value = 72.97991
$detailed.Cells.Item(50, 5).Value = "OFF"

$detailed.Cells.Item(51, 2).Value = 64.89
$detailed.Cells.Item(51, 5).Value = "OFF"

$detailed.Cells.Item(53, 2).Value = 56.98

$detailed.Cells.Item(54, 2).Value = 56.98

$detailed.Cells.Item(55, 2).Value = 52.64259
$detailed.Cells.Item(55, 5).Value = "ON"

$detailed.Cells.Item(56, 2).Value = 49.99118
$detailed.Cells.Item(56, 5).Value = "ON"

$detailed.Cells.Item(57, 2).Value = 36.07
$detailed.Cells.Item(57, 5).Value = "ON"

$detailed.Cells.Item(58, 2).Value = 50.60094
$detailed.Cells.Item(58, 5).Value = "ON"

$detailed.Cells.Item(59, 5).Value = "ON"

$detailed.Cells.Item(60, 2).Value = 60.48199

$detailed.Cells.Item(63, 2).Value = 78.55543

$detailed.Cells.Item(65, 2).Value = 52.38598

$detailed.Cells.Item(66, 2).Value = 56.81604

$detailed.Cells.Item(67, 2).Value = 47.46014
$detailed.Cells.Item(67, 5).Value = "OFF"

$detailed.Cells.Item(68, 2).Value = 36.07

$detailed.Cells.Item(69, 2).Value = 33.0391

$detailed.Cells.Item(70, 2).Value = 20.06478

$detailed.Cells.Item(72, 2).Value = 22.07

$detailed.Cells.Item(73, 2).Value = 0.51

$detailed.Cells.Item(74, 2).Value = 8.038040000000001

$detailed.Cells.Item(75, 2).Value = 0

$detailed.Cells.Item(76, 2).Value = -0.99369

$detailed.Cells.Item(77, 2).Value = -2.60914

$detailed.Cells.Item(78, 2).Value = 0.00831

$detailed.Cells.Item(79, 2).Value = -11.64094

$detailed.Cells.Item(80, 2).Value = 0.51

$detailed.Cells.Item(81, 2).Value = 36.0601

$detailed.Cells.Item(82, 2).Value = 44.9856

$detailed.Cells.Item(84, 2).Value = 47.91478
$detailed.Cells.Item(84, 5).Value = "OFF"

$detailed.Cells.Item(85, 2).Value = 57.09
$detailed.Cells.Item(85, 5).Value = "OFF"

$detailed.Cells.Item(86, 2).Value = 57.09

$detailed.Cells.Item(88, 2).Value = 58.1749

$detailed.Cells.Item(89, 2).Value = 64.89

$detailed.Cells.Item(90, 2).Value = 64.89

$detailed.Cells.Item(91, 2).Value = 64.89

$detailed.Cells.Item(92, 2).Value = 64.89

$detailed.Cells.Item(93, 2).Value = 57.60652

$detailed.Cells.Item(94, 2).Value = 58.42788

$detailed.Cells.Item(95, 2).Value = 57.09

$detailed.Cells.Item(96, 2).Value = 56.98

$detailed.Cells.Item(97, 2).Value = 56.98
